$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "not_auth_1"
$ws.Range("B2").Value = "ТС-А1, ТС-А2, ТС-А3, ТС-А4"

$ws.Range("A3").Value = "not_auth_2"
$ws.Range("B3").Value = "ТС-Р1, ТС-Р2, ТС-Р3"

$ws.Range("A4").Value = "not_auth_3"
$ws.Range("B4").Value = "ТС-ОБ1, ТС-ОБ2"

$ws.Range("A5").Value = "not_auth_4"
$ws.Range("B5").Value = "ТС-ПР1, ТС-ПР2, ТС-ПР3"

$ws.Range("A6").Value = "not_auth_5"
$ws.Range("B6").Value = "6.5.1 ТС-ФИ1 "

$ws.Range("B11").Select()
